# Estadisticos Segundo Parcial 26 Mayo
# Update statistics (Aprobados, Reprobados, percentages, Promedio, Blancos)
# for rows 8-10 on the "2o Parcial" and "Final" sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "2o Parcial" ----
$ws2 = $wb.Worksheets.Item("2o Parcial")

# Row 8 - Saucedo Rivalcoba Liliana Guadalupe / TEMAS DE ADMINISTRACION / 6ARHM
$ws2.Range("E8").Value = 39
$ws2.Range("F8").Value = 0
$ws2.Range("G8").Value = 100
$ws2.Range("H8").Value = 0
$ws2.Range("I8").Value = 9.300000000000001
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0

# Row 9 - Totales Docente
$ws2.Range("E9").Value = 39
$ws2.Range("F9").Value = 0
$ws2.Range("G9").Value = 100
$ws2.Range("H9").Value = 0
$ws2.Range("I9").Value = 9.300000000000001
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 0

# Row 10 - Totales Generales
$ws2.Range("E10").Value = 124
$ws2.Range("F10").Value = 2
$ws2.Range("G10").Value = 98.40000000000001
$ws2.Range("H10").Value = 1.6
$ws2.Range("I10").Value = 8.800000000000001
$ws2.Range("J10").Value = 0
$ws2.Range("K10").Value = 0

# ---- Sheet "Final" ----
$ws3 = $wb.Worksheets.Item("Final")

# Row 8 - Saucedo Rivalcoba Liliana Guadalupe / TEMAS DE ADMINISTRACION / 6ARHM
$ws3.Range("E8").Value = 39
$ws3.Range("F8").Value = 0
$ws3.Range("G8").Value = 100
$ws3.Range("H8").Value = 0
$ws3.Range("I8").Value = 9.300000000000001

# Row 9 - Totales Docente
$ws3.Range("E9").Value = 39
$ws3.Range("F9").Value = 0
$ws3.Range("G9").Value = 100
$ws3.Range("H9").Value = 0
$ws3.Range("I9").Value = 9.300000000000001

# Row 10 - Totales Generales
$ws3.Range("E10").Value = 124
$ws3.Range("F10").Value = 2
$ws3.Range("G10").Value = 98.40000000000001
$ws3.Range("H10").Value = 1.6
